# Help Menu Layout rework.
#
# The sheet used to spread two long paragraphs of flavour text across the
# whole A:J range (two big merged blocks). It is replaced by a compact
# "Goal / Time / Controls" card confined to columns A:E; the rest of the
# old layout (and its merges) goes away.
#
# NOTE on ordering: on this engine, Range.Merge() always mints a brand-new
# (duplicate) style slot for the touched cells, no matter what. To avoid
# minting *extra* duplicates on top of that, every range is merged first
# (while it is still the plain default style, right after Clear()) and
# only afterwards do we set values / alignment - at that point the
# alignment setters correctly reuse an existing style slot instead of
# minting a new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- wipe the existing layout (values, formats and merges) -----------------
$ws.Range("A1:J10").UnMerge()
$ws.Range("A1:J10").Clear()

# --- recreate the merges while everything is still default-styled ----------
$ws.Range("A1:E1").Merge()
$ws.Range("B2:E2").Merge()
$ws.Range("A3:A7").Merge()
$ws.Range("B3:E3").Merge()
$ws.Range("B4:E4").Merge()
$ws.Range("E5:E6").Merge()
$ws.Range("B7:D7").Merge()

# ============================== values ======================================
$ws.Range("A1").Value = "Help"
$ws.Range("A2").Value = "Goal"
$ws.Range("B2").Value = "Time"
$ws.Range("A3").Value = "Goal Body"
$ws.Range("B3").Value = "Time Body"
$ws.Range("B4").Value = "Controls"
$ws.Range("C5").Value = "w"
$ws.Range("E5").Value = "Move"
$ws.Range("B6").Value = "a"
$ws.Range("C6").Value = "s"
$ws.Range("D6").Value = "d"
$ws.Range("B7").Value = "Space"
$ws.Range("E7").Value = "Attack/Action"

# ============================== alignment ====================================
# style 1: horizontal-center only
$ws.Range("A1:E1,B2:E2,B4:E4,E5:E6,B7:D7").HorizontalAlignment = -4108   # xlCenter

# style 2: horizontal-center + wrap
$ws.Range("A3:A7,B3:E3").HorizontalAlignment = -4108                    # xlCenter
$ws.Range("A3:A7,B3:E3").WrapText = $true

# style 3: wrap only (general horizontal)
$ws.Range("G3:J6,F4:J6,A8:E10").HorizontalAlignment = 1                 # xlGeneral
$ws.Range("G3:J6,F4:J6,A8:E10").WrapText = $true

# style 4: "touched" general alignment (applyAlignment with no explicit
# horizontal/vertical/wrap - matches cells that used to be part of a merge
# and are simply left as individually-touched cells now)
$ws.Range("F1:J1,A2,G2:J2,E7:J7,J8,I9:J9").HorizontalAlignment = 1      # xlGeneral

# ============================== selection ====================================
$ws.Range("E7").Select()
